$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $s
}

$ws.Range("D2").Value = '56.707.16'
$ws.Range("E2").Value = '  +3.14%  '
$ws.Range("D3").Value = '2.467.48'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  -0.26%  '
Set-TextValue "D5" '489.41'
$ws.Range("E5").Value = '  +1.93%  '
Set-TextValue "D6" '151.41'
$ws.Range("E6").Value = '  +8.09%  '
Set-TextValue "D7" '0.997'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").Value = '2.474.76'
$ws.Range("E9").Value = '  +0.54%  '
Set-TextValue "D10" '0.0993'
$ws.Range("E10").Value = '  +3.01%  '
Set-TextValue "D11" '5.70'
$ws.Range("E11").Value = '  +3.92%  '
$ws.Range("E12").Value = '  +2.91%  '
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '2.903.24'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").Value = '57.049.92'
$ws.Range("E15").Value = '  +3.17%  '
Set-TextValue "D16" '20.96'
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").Value = '2.474.27'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  +4.49%  '
Set-TextValue "D20" '10.13'
$ws.Range("E20").Value = '  +2.58%  '
Set-TextValue "D21" '319.84'
$ws.Range("E21").Value = '  +1.60%  '
Set-TextValue "D22" '0.998'
Set-TextValue "D23" '5.84'
$ws.Range("E23").Value = '  +3.25%  '
Set-TextValue "D24" '57.88'
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  -0.17%  '
Set-TextValue "D27" '0.162'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").Value = '2.587.85'
$ws.Range("E28").Value = '  +1.87%  '
Set-TextValue "D29" '7.55'
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("D30").Value = '0.0₃0802'
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("E31").Value = '  +0.22%  '
Set-TextValue "D32" '150.87'
$ws.Range("E32").Value = '  +1.53%  '
Set-TextValue "D33" '18.24'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E36").Value = '  +2.25%  '
Set-TextValue "D37" '0.885'
$ws.Range("E37").Value = '  +5.11%  '
$ws.Range("E38").Value = '  +4.70%  '
Set-TextValue "D39" '34.06'
$ws.Range("E39").Value = '  +1.60%  '
$ws.Range("E40").Value = '  +7.23%  '
$ws.Range("E41").Value = '  +2.52%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D42" '0.0557'
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D43" '0.996'
$ws.Range("E43").Value = '  +0.38%  '
Set-TextValue "D44" '0.607'
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("E45").Value = '  +5.50%  '
Set-TextValue "D46" '4.77'
$ws.Range("E46").Value = '  +2.85%  '
Set-TextValue "D47" '261.48'
$ws.Range("E47").Value = '  +3.05%  '
Set-TextValue "D48" '10.21'
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E49").Value = '  +2.54%  '
Set-TextValue "D50" '17.71'
$ws.Range("E50").Value = '  +3.09%  '
$ws.Range("D51").Value = '1.853.45'
$ws.Range("E51").Value = '  -3.43%  '
